$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers I1 ("I0") and J1 ("IF"), styled like the other header cells
# (bold, centered, bordered) by copying H1's format onto the new cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-78 for columns I (I0) and J (IF)
$iValues = @(9, 9, 10, 9, 8, 8, 8, 9, 8, 9, 9, 9, 9, 8, 8, 9, 9, 8, 9, 8, 7, 6, 7, 6, 5, 8, 7, 8, 9, 8, 7, 7, 7, 9, 8, 8, 4, 6, 6, 6, 6, 7, 8, 9, 7, 7, 8, 8, 9, 9, 4, 6, 8, 9, 8, 8, 9, 9, 9, 7, 9, 8, 8, 8, 9, 9, 9, 9, 9, 9, 9, 9, 8, 8, 9, 4, 3)
$jValues = @(9, 9, 11, 9, 8, 9, 9, 9, 9, 9, 9, 9, 9, 9, 8, 9, 9, 9, 9, 8, 8, 7, 7, 7, 6, 8, 7, 8, 9, 8, 7, 8, 7, 9, 8, 8, 4, 7, 7, 6, 6, 7, 8, 9, 8, 8, 8, 8, 9, 9, 6, 7, 8, 9, 8, 8, 9, 9, 9, 7, 9, 9, 9, 8, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 4, 3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
